# Update the Maltaspor player roster table (A1:C19) on the active sheet.
# The player list, position(s), and team have been refreshed; row count
# stays the same (18 data rows + 1 header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Jordan Hawkins", "SG", "New Orleans Pelicans"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Ausar Thompson", "SF,PF", "Detroit Pistons"),
    @("Aaron Nesmith", "SF,PF", "Indiana Pacers"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Brandon Clarke", "PF,C", "Memphis Grizzlies"),
    @("Zach Edey", "C", "Memphis Grizzlies"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Brandon Williams", "PG", "Dallas Mavericks"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Andrew Wiggins", "SF,PF", "Miami Heat"),
    @("Anthony Davis", "PF,C", "Dallas Mavericks"),
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
